$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.701.27"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.599.67"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.57"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "1.824.84"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.600.38"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").Value = "26.677.79"
$ws.Range("D18").Value = "0.0₃0758"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "1.292.12"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.620"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +15.99%  "
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.41"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.21"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "1.737.06"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.57%  "
